# Casos de uso.xlsx - "Agregar boleto de colectivo (falta testearlo)"
#
# Updates the use-case table on Hoja1:
#  - Row 11 (Nro 10): clarifies the description of agregarBoleto (colectivo).
#  - Rows 18-19 (Nro 17-18): the "CalculoDeBeneficios" actor is renamed to
#    "Funciones" and calcularRedSube / calcularTarifaSocial are reworked to
#    operate on a "nuevoBoleto" parameter instead of returning a separate
#    monto value.
#  - Three new rows are appended for the MovimientoAlta actor (agregarBoleto,
#    traerBoletosRedSube) and a new traerTarifaSocial use case under TarjetaABM.
#  - Minor view cosmetics: wider column D, zoom 85%, and the active
#    selection moved to D19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11 (Nro 10): AdminDeLectoras / agregarBoleto ---------------------
$ws.Range("D11").Value = "Crea un objeto boleto con tramo de colectivo y le resta el saldo a la tarjeta."

# --- Row 18 (Nro 17): Funciones / calcularRedSube --------------------------
$ws.Range("B18").Value = "Funciones"
$ws.Range("D18").Value = "Al nuevo boleto le cambia el monto (si corresponde) revisando los boletos de las ultimas 2 horas."
$ws.Range("E18").Value = "boletos:List<Boleto>, nuevoBoleto: Boleto "
$ws.Range("F18").ClearContents()

# --- Row 19 (Nro 18): Funciones / calcularTarifaSocial ---------------------
$ws.Range("B19").Value = "Funciones"
$ws.Range("D19").Value = "Al nuevo boleto le cambia el monto (si corresponde) descontando el porcentaje de la tarifa social."
$ws.Range("E19").Value = "nuevoBoleto: Boleto"
# F19 (montoTarifaSocial: float) is unchanged.

# --- New row 26 (Nro 25): MovimientoAlta / agregarBoleto -------------------
$ws.Range("A26").Value = 25
$ws.Range("B26").Value = "MovimientoAlta"
$ws.Range("C26").Value = "agregarBoleto"
$ws.Range("D26").Value = "agrega un boleto a la base de datos"
$ws.Range("E26").Value = "boleto:Boleto"

# --- New row 27 (Nro 26): MovimientoAlta / traerBoletosRedSube -------------
$ws.Range("A27").Value = 26
$ws.Range("B27").Value = "MovimientoAlta"
$ws.Range("C27").Value = "traerBoletosRedSube"

# --- New row 28 (Nro 27): TarjetaABM / traerTarifaSocial --------------------
$ws.Range("A28").Value = 27
$ws.Range("B28").Value = "TarjetaABM"
$ws.Range("C28").Value = "traerTarifaSocial"

# --- View cosmetics ----------------------------------------------------
# Column D grows to fit the longer descriptions now in it.
$ws.Columns.Item(4).ColumnWidth = 95

# Zoom to 85% and move the selection to D19, matching the saved view state.
$excel.ActiveWindow.Zoom = 85
$ws.Range("D19").Select()
